# Update "想去人数" (want-to-go count) values in F column across sheets
# to reflect newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5605
$ws1.Range("F4").Value = 79
$ws1.Range("F5").Value = 6
$ws1.Range("F6").Value = 928
$ws1.Range("F7").Value = 149
$ws1.Range("F8").Value = 2511
$ws1.Range("F9").Value = 82
$ws1.Range("F10").Value = 129
$ws1.Range("F11").Value = 5
$ws1.Range("F12").Value = 76
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 2347
$ws1.Range("F15").Value = 319

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 103
$ws2.Range("F3").Value = 3

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5605
$ws4.Range("F4").Value = 103
$ws4.Range("F5").Value = 79
$ws4.Range("F6").Value = 6
$ws4.Range("F7").Value = 3
$ws4.Range("F8").Value = 928
$ws4.Range("F9").Value = 149
$ws4.Range("F10").Value = 2511
$ws4.Range("F11").Value = 82
$ws4.Range("F12").Value = 129
$ws4.Range("F13").Value = 5
$ws4.Range("F15").Value = 76
$ws4.Range("F16").Value = 14
$ws4.Range("F17").Value = 2347
$ws4.Range("F18").Value = 319
